# Natmi following Dr Hou advice: expand Slamf7-Slamf7 LR pairs to M2/sCs cluster combinations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M2"
$ws.Range("B2").Value = "Slamf7"
$ws.Range("C2").Value = "Slamf7"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 94.418046
$ws.Range("H2").Value = 283.254138
$ws.Range("I2").Value = 0.9930721842318497
$ws.Range("J2").Value = 0.9930721842318498
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 94.418046
$ws.Range("N2").Value = 283.254138
$ws.Range("O2").Value = 0.9930721842318497
$ws.Range("P2").Value = 0.9930721842318498
$ws.Range("Q2").Value = 8914.767410458117
$ws.Range("R2").Value = 80232.90669412306
$ws.Range("S2").Value = 0.9861923630950169
$ws.Range("T2").Value = 0.9861923630950171

# Row 3
$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Slamf7"
$ws.Range("C3").Value = "Slamf7"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 94.418046
$ws.Range("H3").Value = 283.254138
$ws.Range("I3").Value = 0.9930721842318497
$ws.Range("J3").Value = 0.9930721842318498
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.658674
$ws.Range("N3").Value = 1.976022
$ws.Range("O3").Value = 0.006927815768150184
$ws.Range("P3").Value = 0.006927815768150184
$ws.Range("Q3").Value = 62.190712031004
$ws.Range("R3").Value = 559.716408279036
$ws.Range("S3").Value = 0.006879821136832753
$ws.Range("T3").Value = 0.006879821136832754

# Row 4
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Slamf7"
$ws.Range("C4").Value = "Slamf7"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.658674
$ws.Range("H4").Value = 1.976022
$ws.Range("I4").Value = 0.006927815768150184
$ws.Range("J4").Value = 0.006927815768150184
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 94.418046
$ws.Range("N4").Value = 283.254138
$ws.Range("O4").Value = 0.9930721842318497
$ws.Range("P4").Value = 0.9930721842318498
$ws.Range("Q4").Value = 62.190712031004
$ws.Range("R4").Value = 559.716408279036
$ws.Range("S4").Value = 0.006879821136832753
$ws.Range("T4").Value = 0.006879821136832754

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Slamf7"
$ws.Range("C5").Value = "Slamf7"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.658674
$ws.Range("H5").Value = 1.976022
$ws.Range("I5").Value = 0.006927815768150184
$ws.Range("J5").Value = 0.006927815768150184
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.658674
$ws.Range("N5").Value = 1.976022
$ws.Range("O5").Value = 0.006927815768150184
$ws.Range("P5").Value = 0.006927815768150184
$ws.Range("Q5").Value = 0.433851438276
$ws.Range("R5").Value = 3.904662944484
$ws.Range("S5").Value = 0.00004799463131743032
$ws.Range("T5").Value = 0.00004799463131743032
